$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 65.45238095238095
$ws.Range("H2").Value = 17.52119047619048
$ws.Range("I2").Value = 167.8352380952381
$ws.Range("J2").Value = 7147.117142857142
$ws.Range("G3").Value = 90.5
$ws.Range("H3").Value = 24.57071428571428
$ws.Range("I3").Value = 186.4028571428571
$ws.Range("J3").Value = 9895.018571428573
$ws.Range("G4").Value = 96.78571428571429
$ws.Range("H4").Value = 28.86642857142857
$ws.Range("I4").Value = 197.6928571428571
$ws.Range("J4").Value = 12880.15142857143
$ws.Range("G5").Value = 68
$ws.Range("H5").Value = 17.10777777777778
$ws.Range("I5").Value = 124.8822222222222
$ws.Range("J5").Value = 6191.664444444445
$ws.Range("G6").Value = 133.7222222222222
$ws.Range("H6").Value = 33.53666666666666
$ws.Range("I6").Value = 128.4011111111111
$ws.Range("J6").Value = 13032.63444444444
$ws.Range("G7").Value = 111.3125
$ws.Range("H7").Value = 30.713125
$ws.Range("I7").Value = 155.031875
$ws.Range("J7").Value = 13145.564375
$ws.Range("G8").Value = 134.7142857142857
$ws.Range("H8").Value = 38.19928571428571
$ws.Range("I8").Value = 161.6435714285714
$ws.Range("J8").Value = 16368.36285714286
$ws.Range("G9").Value = 278.625
$ws.Range("H9").Value = 66.910625
$ws.Range("I9").Value = 117.86875
$ws.Range("J9").Value = 24165.876875
$ws.Range("G10").Value = 191.3125
$ws.Range("H10").Value = 52.52125
$ws.Range("I10").Value = 154.361875
$ws.Range("J10").Value = 20433.32125
$ws.Range("G11").Value = 165.9375
$ws.Range("H11").Value = 43.170625
$ws.Range("I11").Value = 125.0475
$ws.Range("J11").Value = 15489.34625
$ws.Range("G12").Value = 249
$ws.Range("H12").Value = 58.83428571428572
$ws.Range("I12").Value = 84.25142857142858
$ws.Range("J12").Value = 18264.09428571429
$ws.Range("G13").Value = 198.5625
$ws.Range("H13").Value = 52.12125
$ws.Range("I13").Value = 125.6175
$ws.Range("J13").Value = 19100.43375
$ws.Range("G14").Value = 106.5
$ws.Range("H14").Value = 26.948125
$ws.Range("I14").Value = 121.74625
$ws.Range("J14").Value = 9708.01125
$ws.Range("G15").Value = 135.9166666666667
$ws.Range("H15").Value = 27.24541666666667
$ws.Range("I15").Value = 63.51666666666667
$ws.Range("J15").Value = 7768.971666666666
$ws.Range("G16").Value = 149.8
$ws.Range("H16").Value = 28.383
$ws.Range("I16").Value = 46.726
$ws.Range("J16").Value = 7865.425999999999
$ws.Range("G17").Value = 79.4375
$ws.Range("H17").Value = 20.95125
$ws.Range("I17").Value = 150.763125
$ws.Range("J17").Value = 8334.363125
$ws.Range("G18").Value = 144.4444444444445
$ws.Range("H18").Value = 38.96777777777778
$ws.Range("I18").Value = 125.4611111111111
$ws.Range("J18").Value = 15320.04833333333
$ws.Range("G19").Value = 140.1111111111111
$ws.Range("H19").Value = 34.17555555555555
$ws.Range("I19").Value = 108.1377777777778
$ws.Range("J19").Value = 11662.68888888889
